$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.828.09'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").Value = '3.366.35'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.12'
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.07'
$ws.Range("E6").Value = '  -1.98%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.472'
$ws.Range("E8").Value = '  -0.48%  '
$ws.Range("E9").Value = '  +2.55%  '
$ws.Range("E10").Value = '  -1.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.382'
$ws.Range("E11").Value = '  -3.23%  '
$ws.Range("D12").Value = '3.938.46'
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("E13").Value = '  +0.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.13'
$ws.Range("E14").Value = '  -1.24%  '
$ws.Range("D15").Value = '3.366.92'
$ws.Range("E15").Value = '  -0.46%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000168'
$ws.Range("E16").Value = '  -1.54%  '
$ws.Range("D17").Value = '60.936.97'
$ws.Range("E17").Value = '  +0.16%  '
$ws.Range("E18").Value = '  -1.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.51'
$ws.Range("E19").Value = '  -3.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.91'
$ws.Range("E20").Value = '  -0.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '383.64'
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.01'
$ws.Range("E22").Value = '  +1.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.550'
$ws.Range("E23").Value = '  -1.48%  '
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("E25").Value = '  -5.72%  '
$ws.Range("E26").Value = '  +6.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.09%  '
$ws.Range("E28").Value = '  -3.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.95'
$ws.Range("E29").Value = '  -0.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.11'
$ws.Range("E30").Value = '  -1.63%  '
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("E32").Value = '  -6.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.02'
$ws.Range("E33").Value = '  -2.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.82'
$ws.Range("E34").Value = '  -1.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '167.16'
$ws.Range("E35").Value = '  +0.35%  '
$ws.Range("E36").Value = '  -1.11%  '
$ws.Range("D37").Value = '3.402.15'
$ws.Range("E37").Value = '  -0.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.44'
$ws.Range("E38").Value = '  -2.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0754'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '25.75'
$ws.Range("E40").Value = '  -7.64%  '
$ws.Range("E41").Value = '  -1.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.34'
$ws.Range("E42").Value = '  -1.43%  '
$ws.Range("E43").Value = '  -2.72%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.11'
$ws.Range("E44").Value = '  -1.61%  '
$ws.Range("D45").Value = '2.437.43'
$ws.Range("E45").Value = '  -2.13%  '
$ws.Range("B46").Value = 'Cosmos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.62'
$ws.Range("E46").Value = '  -2.64%  '
$ws.Range("B47").Value = 'FirstDigitalUSD'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.05'
$ws.Range("E48").Value = '  -6.48%  '
$ws.Range("E49").Value = '  -5.17%  '
$ws.Range("E50").Value = '  -5.01%  '
$ws.Range("E51").Value = '  -2.90%  '
